$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'302.71"
$ws.Range("E2").Value = "'1.97%"
$ws.Range("D3").Value = "'31.75"
$ws.Range("E3").Value = "'0.36%"
$ws.Range("D4").Value = "'5.165"
$ws.Range("D5").Value = "'0.07816"
$ws.Range("E5").Value = "'4.43%"
$ws.Range("D6").Value = "'2.312"
$ws.Range("E6").Value = "'34.94%"
$ws.Range("D7").Value = "'7.944"
$ws.Range("E7").Value = "'2.61%"
$ws.Range("B8").Value = "'MXToken"
$ws.Range("C8").Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.9071"
$ws.Range("E8").Value = "'-2.59%"
$ws.Range("B9").Value = "'WazirX"
$ws.Range("C9").Value = "'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").Value = "'0.1738"
$ws.Range("E9").Value = "'2.05%"
$ws.Range("B10").Value = "'LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.07367"
$ws.Range("E10").Value = "'1.74%"
$ws.Range("B11").Value = "'MandalaExchangeToken"
$ws.Range("C11").Value = "'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.08136"
$ws.Range("E11").Value = "'2.58%"
$ws.Range("B12").Value = "'BitrueCoin"
$ws.Range("C12").Value = "'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.03010"
$ws.Range("E12").Value = "'-0.10%"
$ws.Range("B13").Value = "'BitMartToken"
$ws.Range("C13").Value = "'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09940"
$ws.Range("E13").Value = "'0.38%"
$ws.Range("B14").Value = "'BitForexToken"
$ws.Range("C14").Value = "'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001505"
$ws.Range("E14").Value = "'0.63%"
$ws.Range("B15").Value = "'TigerCash"
$ws.Range("C15").Value = "'https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.006096"
$ws.Range("E15").Value = "'-3.19%"
$ws.Range("B16").Value = "'LEO"
$ws.Range("C16").Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.497"
$ws.Range("E16").Value = "'1.27%"
$ws.Range("B17").Value = "'GateToken"
$ws.Range("C17").Value = "'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "'3.871"
$ws.Range("E17").Value = "'1.85%"
$ws.Range("D19").Value = "'0.3243"
$ws.Range("E19").Value = "'-0.93%"
$ws.Range("D20").Value = "'0.1338"
$ws.Range("E20").Value = "'0.81%"
$ws.Range("D21").Value = "'4.674"
$ws.Range("E21").Value = "'2.29%"
$ws.Range("D22").Value = "'0.04647"
$ws.Range("E22").Value = "'-0.15%"
$ws.Range("D23").Value = "'0.1566"
$ws.Range("E23").Value = "'0.42%"
$ws.Range("E24").Value = "'3.67%"
$ws.Range("D25").Value = "'0.004525"
$ws.Range("E25").Value = "'2.23%"
$ws.Range("E39").Value = "'8.34%"
$ws.Range("E40").Value = "'2.56%"
$ws.Range("D41").Value = "'0.007298"
$ws.Range("E41").Value = "'3.30%"
$ws.Range("D42").Value = "'0.1361"
$ws.Range("E42").Value = "'2.42%"
$ws.Range("D43").Value = "'0.002249"
$ws.Range("E43").Value = "'9.13%"
$ws.Range("D44").Value = "'0.01074"
$ws.Range("E44").Value = "'-5.42%"
$ws.Range("D45").Value = "'0.00006501"
$ws.Range("E45").Value = "'7.68%"
$ws.Range("E46").Value = "'-0.04%"
$ws.Range("E47").Value = "'15.31%"
$ws.Range("E49").Value = "'-0.04%"
$ws.Range("E50").Value = "'0.03%"
